$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: the phone number in A19 was previously stored as text; it is now
# corrected to a real number (matches the numeric phone cells above it).
$ws.Range("A19").Value = 71277620

# Row 20: new redemption event row for phone 71277620 / 76 points, appended
# below the existing data. The phone number here is entered as raw text
# (same shape the original A19 had before its own correction), so force a
# text entry via a leading apostrophe and then drop the resulting
# "stored as text" formatting so no extra style sticks to the cell.
$ws.Range("A20").Value = "'71277620"
$ws.Range("A20").ClearFormats()

$ws.Range("B20").Value = 76
$ws.Range("C20").Value = "2025-08-18T17:28:56"
